# rename db to db_proxy
# - ServerID for RedisServer_1 goes from "001" to "1"
# - Auth value renamed from "indieleague" to "pwnsky_squick" (both RedisServer_1 and the
#   new MysqlServer_1 row use it)
# - a new MysqlServer_1 row is appended (row 12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 11 (RedisServer_1): ServerID 001 -> 1, Auth indieleague -> pwnsky_squick ---
$ws.Cells.Item(11, 2).Value = "1"

# Match the saved file's dedup of the (identical) trailing cell style so F11 lines up
# on the same xf as the rest of the row instead of its own now-redundant copy.
$ws.Cells.Item(11, 5).Copy()
$ws.Cells.Item(11, 6).PasteSpecial(-4122)
$ws.Cells.Item(11, 6).Value = "pwnsky_squick"

# --- row 12 (new): MysqlServer_1 ---
$ws.Cells.Item(12, 1).Value = "MysqlServer_1"
$ws.Cells.Item(12, 2).Value = 2
$ws.Cells.Item(12, 5).Value = 10086
$ws.Cells.Item(12, 6).Value = "pwnsky_squick"

# IP / PublicIP columns reuse the same style as row 11's (numFmtId 49 / left-center align)
$ws.Cells.Item(11, 3).Copy()
$ws.Range($ws.Cells.Item(12, 3), $ws.Cells.Item(12, 4)).PasteSpecial(-4122)
$ws.Cells.Item(12, 3).Value = "1.14.123.62"
$ws.Cells.Item(12, 4).Value = "1.14.123.62"

# --- cosmetic bits that came along with the resave ---
$ws.Range("G12").Select()

$ws.Rows.Item(10).RowHeight = 14.25
